# ---------------------------------------------------------------------------
# This script reproduces (as closely as the COM surface allows) the edit
# described by the commit: "EXCEL from last commit"
#
#  1) Workbook-level absolute-path hint (x15ac:absPath) changed
#     (not exposed via the Excel object model - informational metadata only,
#      left untouched).
#  2) A brand-new worksheet "Sheet3" is appended after "Sheet2" and becomes
#     the active/selected sheet.
#  3) That new sheet is populated with PGN/J1939-style parameter listings
#     for three PGNs (TSC1, EEC1, DM1), including helper CONCAT() formulas
#     for the TSC1 block.
#  4) Column widths on the new sheet are best-fit to their content.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet as the very last tab (after "Sheet2") --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
# Excel will name it "Sheet3" automatically (next free default sheet name).

# --- 2. Block 1 : TSC1 parameters, rows 3-12 --------------------------------
# Cell values are entered in the same order the shared-string table was
# originally built in (row 4 before row 3, then top-to-bottom), so that the
# generated shared-string indices line up with the source workbook.
$new.Range("C4").Value  = "EngRqedSpeedCtrlConditions"
$new.Range("C3").Value  = "EngOverrideCtrlMode"
$new.Range("C5").Value  = "OverrideCtrlModePriority"
$new.Range("C6").Value  = "EngRqedSpeed_SpeedLimit"
$new.Range("C7").Value  = "EngRqedTorque_TorqueLimit"
$new.Range("C8").Value  = "TransmissionRate"
$new.Range("C9").Value  = "ControlPurpose"
$new.Range("C10").Value = "EngineRequestedTorqueHiRes"
$new.Range("C11").Value = "MessageCounter"
$new.Range("C12").Value = "MessageChecksum"

$new.Range("E3").Value  = "USINT"
$new.Range("E4").Value  = "USINT"
$new.Range("E5").Value  = "USINT"
$new.Range("E6").Value  = "REAL"
$new.Range("E7").Value  = "INT"
$new.Range("E8").Value  = "USINT"
$new.Range("E9").Value  = "USINT"
$new.Range("E10").Value = "REAL"
$new.Range("E11").Value = "USINT"
$new.Range("E12").Value = "USINT"

$new.Range("B3").Value  = "TSC1"
$new.Range("B4").Value  = "TSC1"
$new.Range("B5").Value  = "TSC1"
$new.Range("B6").Value  = "TSC1"
$new.Range("B7").Value  = "TSC1"
$new.Range("B8").Value  = "TSC1"
$new.Range("B9").Value  = "TSC1"
$new.Range("B10").Value = "TSC1"
$new.Range("B11").Value = "TSC1"
$new.Range("B12").Value = "TSC1"

$new.Range("H3").Formula  = "=CONCAT(B3,C3,"":"",E3,"";"")"
$new.Range("H4").Formula  = "=CONCAT(B4,C4,"":"",E4,"";"")"
$new.Range("H5").Formula  = "=CONCAT(B5,C5,"":"",E5,"";"")"
$new.Range("H6").Formula  = "=CONCAT(B6,C6,"":"",E6,"";"")"
$new.Range("H7").Formula  = "=CONCAT(B7,C7,"":"",E7,"";"")"
$new.Range("H8").Formula  = "=CONCAT(B8,C8,"":"",E8,"";"")"
$new.Range("H9").Formula  = "=CONCAT(B9,C9,"":"",E9,"";"")"
$new.Range("H10").Formula = "=CONCAT(B10,C10,"":"",E10,"";"")"
$new.Range("H11").Formula = "=CONCAT(B11,C11,"":"",E11,"";"")"
$new.Range("H12").Formula = "=CONCAT(B12,C12,"":"",E12,"";"")"

# --- 3. Block 2 : EEC1 parameters, rows 15-22 -------------------------------
$new.Range("C15").Value = "EngTorqueMode"
$new.Range("C16").Value = "ActlEngPrcntTrqueHighResolution"
$new.Range("C17").Value = "DriversDemandEngPercentTorque"
$new.Range("C18").Value = "ActualEngPercentTorque"
$new.Range("C19").Value = "EngSpeed"
$new.Range("C20").Value = "SrcAddrssOfCtrllngDvcFrEngCntrl"
$new.Range("C21").Value = "EngStarterMode"
$new.Range("C22").Value = "EngDemandPercentTorque"

$new.Range("E15").Value = "USINT"
$new.Range("E16").Value = "REAL"
$new.Range("E17").Value = "INT"
$new.Range("E18").Value = "INT"
$new.Range("E19").Value = "REAL"
$new.Range("E20").Value = "USINT"
$new.Range("E21").Value = "USINT"
$new.Range("E22").Value = "INT"

$new.Range("B15").Value = "EEC1"
$new.Range("B16").Value = "EEC1"
$new.Range("B17").Value = "EEC1"
$new.Range("B18").Value = "EEC1"
$new.Range("B19").Value = "EEC1"
$new.Range("B20").Value = "EEC1"
$new.Range("B21").Value = "EEC1"
$new.Range("B22").Value = "EEC1"

# --- 4. Block 3 : DM1 lamp-status parameters, rows 26-33 --------------------
$new.Range("C26").Value = "ProtectLampStatus"
$new.Range("C27").Value = "AmberWarningLampStatus"
$new.Range("C28").Value = "RedStopLampState"
$new.Range("C29").Value = "MalfunctionIndicatorLampStatus"
$new.Range("C30").Value = "FlashProtectLamp"
$new.Range("C31").Value = "FlashAmberWarningLamp"
$new.Range("C32").Value = "FlashRedStopLamp"
$new.Range("C33").Value = "FlashMalfuncIndicatorLamp"

$new.Range("E26").Value = "USINT"
$new.Range("E27").Value = "USINT"
$new.Range("E28").Value = "USINT"
$new.Range("E29").Value = "USINT"
$new.Range("E30").Value = "USINT"
$new.Range("E31").Value = "USINT"
$new.Range("E32").Value = "USINT"
$new.Range("E33").Value = "USINT"

$new.Range("B26").Value = "DM1"
$new.Range("B27").Value = "DM1"
$new.Range("B28").Value = "DM1"
$new.Range("B29").Value = "DM1"
$new.Range("B30").Value = "DM1"
$new.Range("B31").Value = "DM1"
$new.Range("B32").Value = "DM1"
$new.Range("B33").Value = "DM1"

# --- 5. Block 4 : DM1 diagnostic trouble codes, rows 34-38 ------------------
$new.Range("C34").Value = "DTC1"
$new.Range("C35").Value = "DTC2"
$new.Range("C36").Value = "DTC3"
$new.Range("C37").Value = "DTC4"
$new.Range("C38").Value = "DTC5"

$new.Range("E34").Value = "UDINT"
$new.Range("E35").Value = "UDINT"
$new.Range("E36").Value = "UDINT"
$new.Range("E37").Value = "UDINT"
$new.Range("E38").Value = "UDINT"

$new.Range("B34").Value = "DM1"
$new.Range("B35").Value = "DM1"
$new.Range("B36").Value = "DM1"
$new.Range("B37").Value = "DM1"
$new.Range("B38").Value = "DM1"

# --- 6. Cosmetics: best-fit column widths for C and H -----------------------
$new.Columns.Item(3).ColumnWidth = 30.877604166666668
$new.Columns.Item(8).ColumnWidth = 33.451822916666664

# --- 7. Selection / active cell on the new sheet, matching the saved view --
[void]$new.Range("H26").Select()
